# Enhance user authentication module
# Appends a new row (row 46) to each of the 4 worksheets, duplicating the
# last existing row (row 45) and updating the timestamp (column A) plus,
# for the first sheet only, the "actual length" (D) and its decimal
# counterpart (H).

$wb = $excel.ActiveWorkbook

$newTimestamp = 45832.43385416667

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Duplicate row 45 into row 46 (values, formats, styles all carried over)
    $ws.Range("A45:I45").Copy($ws.Range("A46:I46"))

    # New sample date/time
    $ws.Cells.Item(46, 1).Value = $newTimestamp

    if ($i -eq 1) {
        # Sheet "DE_LFT_#1" has a slightly different actual-length reading
        $ws.Cells.Item(46, 4).Value = "0x01,0x64"
        $ws.Cells.Item(46, 8).Value = 356
    }
}
